# Commit: "CAPM, FCF, and Brent Stuff"
# Renames the second sheet "FED" -> "BrentOilPrices" and populates it with
# quarterly Brent crude oil price data (Date / Value columns), reusing the
# ECB sheet's existing date-format cell styles. Also updates the view state:
# BrentOilPrices becomes the active tab (selection B64) while the ECB sheet's
# selection becomes an entire-column selection on column A.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ECB")
$ws2 = $wb.Worksheets.Item("FED")   # renamed to "BrentOilPrices" below

# 1. Rename sheet 2
$ws2.Name = "BrentOilPrices"

# 2. Header row: shared strings "Date" (reused from ECB) and "Value" (new)
$ws2.Cells.Item(1, 1).Value = "Date"
$ws2.Cells.Item(1, 2).Value = "Value"

# 3. Bulk-write the 62 data rows: col A = date serials, col B = prices
$arr = New-Object 'object[,]' 62,2
$arr[0,0] = 43830; $arr[0,1] = 66
$arr[1,0] = 43738; $arr[1,1] = 62.39
$arr[2,0] = 43646; $arr[2,1] = 65.2
$arr[3,0] = 43555; $arr[3,1] = 67.540000000000006
$arr[4,0] = 43465; $arr[4,1] = 59.97
$arr[5,0] = 43373; $arr[5,1] = 79.739999999999995
$arr[6,0] = 43281; $arr[6,1] = 76.459999999999994
$arr[7,0] = 43190; $arr[7,1] = 69.53
$arr[8,0] = 43100; $arr[8,1] = 70.52
$arr[9,0] = 43008; $arr[9,1] = 58.33
$arr[10,0] = 42916; $arr[10,1] = 49.3
$arr[11,0] = 42825; $arr[11,1] = 55.36
$arr[12,0] = 42735; $arr[12,1] = 53.65
$arr[13,0] = 42643; $arr[13,1] = 50.61
$arr[14,0] = 42551; $arr[14,1] = 45.26
$arr[15,0] = 42460; $arr[15,1] = 42.86
$arr[16,0] = 42369; $arr[16,1] = 33.58
$arr[17,0] = 42277; $arr[17,1] = 48.23
$arr[18,0] = 42185; $arr[18,1] = 56.52
$arr[19,0] = 42094; $arr[19,1] = 62.84
$arr[20,0] = 42004; $arr[20,1] = 51.2
$arr[21,0] = 41912; $arr[21,1] = 86.2
$arr[22,0] = 41820; $arr[22,1] = 107.2
$arr[23,0] = 41729; $arr[23,1] = 109.9
$arr[24,0] = 41639; $arr[24,1] = 106.31
$arr[25,0] = 41547; $arr[25,1] = 111.21
$arr[26,0] = 41455; $arr[26,1] = 107.6
$arr[27,0] = 41364; $arr[27,1] = 103.7
$arr[28,0] = 41274; $arr[28,1] = 113.3
$arr[29,0] = 41182; $arr[29,1] = 112.47
$arr[30,0] = 41090; $arr[30,1] = 103.55
$arr[31,0] = 40999; $arr[31,1] = 120.18
$arr[32,0] = 40908; $arr[32,1] = 110.75
$arr[33,0] = 40816; $arr[33,1] = 111.45
$arr[34,0] = 40724; $arr[34,1] = 118.15
$arr[35,0] = 40633; $arr[35,1] = 122
$arr[36,0] = 40543; $arr[36,1] = 95.7
$arr[37,0] = 40451; $arr[37,1] = 79.03
$arr[38,0] = 40359; $arr[38,1] = 74.247
$arr[39,0] = 40268; $arr[39,1] = 81.17
$arr[40,0] = 40178; $arr[40,1] = 72.989999999999995
$arr[41,0] = 40086; $arr[41,1] = 70.19
$arr[42,0] = 39994; $arr[42,1] = 70.11
$arr[43,0] = 39903; $arr[43,1] = 50.32
$arr[44,0] = 39813; $arr[44,1] = 43.62
$arr[45,0] = 39721; $arr[45,1] = 74.53
$arr[46,0] = 39629; $arr[46,1] = 142.03
$arr[47,0] = 39538; $arr[47,1] = 111.36
$arr[48,0] = 39447; $arr[48,1] = 89.07
$arr[49,0] = 39355; $arr[49,1] = 84.6
$arr[50,0] = 39263; $arr[50,1] = 76.290000000000006
$arr[51,0] = 39172; $arr[51,1] = 67.84
$arr[52,0] = 39082; $arr[52,1] = 53.68
$arr[53,0] = 38990; $arr[53,1] = 59.21
$arr[54,0] = 38898; $arr[54,1] = 76.28
$arr[55,0] = 38807; $arr[55,1] = 69.41
$arr[56,0] = 38717; $arr[56,1] = 66.150000000000006
$arr[57,0] = 38625; $arr[57,1] = 57.77
$arr[58,0] = 38533; $arr[58,1] = 57.11
$arr[59,0] = 38442; $arr[59,1] = 50.48
$arr[60,0] = 38352; $arr[60,1] = 44.95
$arr[61,0] = 38260; $arr[61,1] = 50.72
$ws2.Range("A2:B63").Value = $arr

# 4. Reuse ECB's per-row date-format styles for column A (rows 2-63 share the
#    exact same style pattern in both sheets, indices 1 and 2).
for ($r = 2; $r -le 63; $r++) {
    $ws1.Cells.Item($r, 1).Copy() | Out-Null
    $ws2.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$excel.CutCopyMode = $false

# 5. View state: ECB gets an entire-column A selection (no longer the active tab)
$ws1.Activate() | Out-Null
$ws1.Columns.Item(1).Select() | Out-Null

# 6. BrentOilPrices becomes the active tab, with B64 selected
$ws2.Activate() | Out-Null
$ws2.Range("B64").Select() | Out-Null

Write-Host "Edit applied"
